$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 461
$ws.Range("I41").Value = 102.5
$ws.Range("J41").Value = 550.625
$ws.Range("K41").Value = 102.5
$ws.Range("L41").Value = 550.625
$ws.Range("M41").Value = 337.5
$ws.Range("N41").Value = -1430.625
$ws.Range("H129").Value = 100800.31
$ws.Range("I129").Value = 477.875
$ws.Range("J129").Value = 109524
$ws.Range("K129").Value = 1433.625
$ws.Range("L129").Value = 328572
$ws.Range("M129").Value = 3566.375
$ws.Range("N129").Value = -338572
$ws.Range("H132").Value = 3151.5715
$ws.Range("I132").Value = 3404.0454
$ws.Range("J132").Value = 2225.8333
$ws.Range("K132").Value = 10212.1362
$ws.Range("L132").Value = 6677.499899999999
$ws.Range("M132").Value = -7682.136200000001
$ws.Range("N132").Value = -11737.4999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3432.4285
$ws.Range("I32").Value = 3061.2666
$ws.Range("K32").Value = 3061.2666
$ws.Range("M32").Value = -2774.2666
$ws.Range("H61").Value = 2075.361
$ws.Range("I61").Value = 1156.4445
$ws.Range("J61").Value = 4832.1113
$ws.Range("K61").Value = 1156.4445
$ws.Range("L61").Value = 4832.1113
$ws.Range("M61").Value = -944.4445000000001
$ws.Range("N61").Value = -5256.1113
$ws.Range("H88").Value = 201810
$ws.Range("I88").Value = 1665.3334
$ws.Range("J88").Value = 502027
$ws.Range("K88").Value = 1665.3334
$ws.Range("L88").Value = 502027
$ws.Range("M88").Value = -1259.3334
$ws.Range("N88").Value = -502839
$ws.Range("H91").Value = 201810
$ws.Range("I91").Value = 1665.3334
$ws.Range("J91").Value = 502027
$ws.Range("K91").Value = 1665.3334
$ws.Range("L91").Value = 502027
$ws.Range("M91").Value = -261.3334
$ws.Range("N91").Value = -504835
$ws.Range("H122").Value = 3566.6667
$ws.Range("I122").Value = 3750
$ws.Range("K122").Value = 11250
$ws.Range("M122").Value = -8800
$ws.Range("H132").Value = 13345.559
$ws.Range("I132").Value = 1562.1316
$ws.Range("J132").Value = 102899.6
$ws.Range("K132").Value = 4686.3948
$ws.Range("L132").Value = 308698.8
$ws.Range("M132").Value = -2156.3948
$ws.Range("N132").Value = -313758.8
$ws.Range("H136").Value = 2075.361
$ws.Range("I136").Value = 1156.4445
$ws.Range("J136").Value = 4832.1113
$ws.Range("K136").Value = 3469.3335
$ws.Range("L136").Value = 14496.3339
$ws.Range("M136").Value = -919.3335000000002
$ws.Range("N136").Value = -19596.3339

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 336
$ws.Range("I22").Value = 336
$ws.Range("K22").Value = 336
$ws.Range("M22").Value = -163
$ws.Range("H94").Value = 2188.4443
$ws.Range("I94").Value = 1825.4667
$ws.Range("J94").Value = 4003.3333
$ws.Range("K94").Value = 1825.4667
$ws.Range("L94").Value = 4003.3333
$ws.Range("M94").Value = -1374.4667
$ws.Range("N94").Value = -4905.3333
$ws.Range("H105").Value = 2275823.2
$ws.Range("I105").Value = 4290
$ws.Range("J105").Value = 3848423
$ws.Range("K105").Value = 4290
$ws.Range("L105").Value = 3848423
$ws.Range("M105").Value = -2543
$ws.Range("N105").Value = -3851917

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 22222326
$ws.Range("I7").Value = 55555704
$ws.Range("J7").Value = 74
$ws.Range("K7").Value = 55555704
$ws.Range("L7").Value = 74
$ws.Range("M7").Value = -55555591
$ws.Range("N7").Value = -300
$ws.Range("H31").Value = 3116.348
$ws.Range("I31").Value = 1435.875
$ws.Range("J31").Value = 4012.6
$ws.Range("K31").Value = 1435.875
$ws.Range("L31").Value = 4012.6
$ws.Range("M31").Value = -1140.875
$ws.Range("N31").Value = -4602.6
$ws.Range("H34").Value = 3116.348
$ws.Range("I34").Value = 1435.875
$ws.Range("J34").Value = 4012.6
$ws.Range("K34").Value = 1435.875
$ws.Range("L34").Value = 4012.6
$ws.Range("M34").Value = -1233.875
$ws.Range("N34").Value = -4416.6
$ws.Range("H94").Value = 3834.8
$ws.Range("I94").Value = 2540.8
$ws.Range("J94").Value = 5128.8
$ws.Range("K94").Value = 2540.8
$ws.Range("L94").Value = 5128.8
$ws.Range("M94").Value = -2089.8
$ws.Range("N94").Value = -6030.8
$ws.Range("H107").Value = 976.6429000000001
$ws.Range("I107").Value = 836.7778
$ws.Range("K107").Value = 836.7778
$ws.Range("M107").Value = 1083.2222
$ws.Range("H132").Value = 2492.2727
$ws.Range("I132").Value = 1048.8235
$ws.Range("K132").Value = 3146.4705
$ws.Range("M132").Value = -616.4704999999999
$ws.Range("H134").Value = 1211.6
$ws.Range("I134").Value = 908.2222
$ws.Range("K134").Value = 2724.6666
$ws.Range("M134").Value = -189.6666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 118.2
$ws.Range("J12").Value = 125.77778
$ws.Range("L12").Value = 377.33334
$ws.Range("N12").Value = -723.33334
$ws.Range("H95").Value = 5021.6
$ws.Range("I95").Value = 5000
$ws.Range("J95").Value = 5024
$ws.Range("K95").Value = 15000
$ws.Range("L95").Value = 15072
$ws.Range("M95").Value = -12941
$ws.Range("N95").Value = -19190
$ws.Range("H131").Value = 801.0700000000001
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 804.19586
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 2412.58758
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -12492.58758
$ws.Range("H138").Value = 131693.73
$ws.Range("I138").Value = 1227.7693
$ws.Range("J138").Value = 301299.5
$ws.Range("K138").Value = 3683.3079
$ws.Range("L138").Value = 903898.5
$ws.Range("M138").Value = 1456.6921
$ws.Range("N138").Value = -914178.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 59.45
$ws.Range("I2").Value = 44.363636
$ws.Range("J2").Value = 77.888885
$ws.Range("K2").Value = 44.363636
$ws.Range("L2").Value = 77.888885
$ws.Range("M2").Value = 68.636364
$ws.Range("N2").Value = -303.888885
$ws.Range("H132").Value = 68324.75
$ws.Range("I132").Value = 6920
$ws.Range("J132").Value = 170666
$ws.Range("K132").Value = 20760
$ws.Range("L132").Value = 511998
$ws.Range("M132").Value = -18230
$ws.Range("N132").Value = -517058

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2600
$ws.Range("I46").Value = 2000
$ws.Range("K46").Value = 2000
$ws.Range("M46").Value = -1812
$ws.Range("H93").Value = 3500
$ws.Range("I93").Value = 3500
$ws.Range("K93").Value = 3500
$ws.Range("M93").Value = -2252
$ws.Range("H122").Value = 1964709.1
$ws.Range("I122").Value = 2804055.8
$ws.Range("K122").Value = 8412167.399999999
$ws.Range("M122").Value = -8409717.399999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 95464.5
$ws.Range("J138").Value = 95464.5
$ws.Range("L138").Value = 95464.5
$ws.Range("N138").Value = -105744.5
